$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.481.61'
$ws.Range('E2').Value = '  -1.49%  '
$ws.Range('D3').Value = '2.427.01'
$ws.Range('E3').Value = '  -2.31%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'565.75"
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').Value = "'143.54"
$ws.Range('E6').Value = '  -3.69%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -2.01%  '
$ws.Range('D9').Value = '2.425.71'
$ws.Range('E9').Value = '  -2.31%  '
$ws.Range('E10').Value = '  -5.09%  '
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('E12').Value = '  -2.82%  '
$ws.Range('E13').Value = '  -3.15%  '
$ws.Range('D14').Value = "'26.53"
$ws.Range('E14').Value = '  -2.83%  '
$ws.Range('E15').Value = '  -5.86%  '
$ws.Range('D16').Value = '2.862.38'
$ws.Range('E16').Value = '  -2.93%  '
$ws.Range('D17').Value = '62.469.73'
$ws.Range('E17').Value = '  -1.28%  '
$ws.Range('D18').Value = '2.439.41'
$ws.Range('E18').Value = '  -1.93%  '
$ws.Range('E19').Value = '  -4.12%  '
$ws.Range('E20').Value = '  -0.91%  '
$ws.Range('D21').Value = "'324.04"
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('E22').Value = '  -2.46%  '
$ws.Range('E23').Value = '  +8.10%  '
$ws.Range('D24').Value = "'1.00"
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').Value = "'65.01"
$ws.Range('E25').Value = '  -3.84%  '
$ws.Range('D26').Value = "'611.69"
$ws.Range('E26').Value = '  -4.23%  '
$ws.Range('E27').Value = '  +1.44%  '
$ws.Range('D28').Value = '0.0₃0972'
$ws.Range('E28').Value = '  -7.52%  '
$ws.Range('D29').Value = '2.553.63'
$ws.Range('E29').Value = '  -3.87%  '
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('E31').Value = '  -4.01%  '
$ws.Range('E32').Value = '  -4.73%  '
$ws.Range('E33').Value = '  -2.42%  '
$ws.Range('D34').Value = "'0.136"
$ws.Range('E34').Value = '  -5.39%  '
$ws.Range('E35').Value = '  -3.84%  '
$ws.Range('E36').Value = '  -4.93%  '
$ws.Range('D37').Value = "'0.999"
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  -3.17%  '
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('D40').Value = "'147.47"
$ws.Range('E40').Value = '  +0.29%  '
$ws.Range('D41').Value = "'5.24"
$ws.Range('E41').Value = '  -4.85%  '
$ws.Range('E42').Value = '  -6.47%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = "'2.54"
$ws.Range('E43').Value = '  -3.19%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').Value = "'0.999"
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = "'42.10"
$ws.Range('E45').Value = '  +0.74%  '
$ws.Range('D46').Value = "'144.64"
$ws.Range('E46').Value = '  -3.89%  '
$ws.Range('E47').Value = '  -1.46%  '
$ws.Range('D48').Value = "'20.26"
$ws.Range('E48').Value = '  -4.33%  '
$ws.Range('E49').Value = '  -4.38%  '
$ws.Range('E50').Value = '  -2.26%  '
$ws.Range('E51').Value = '  -4.66%  '
